$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("S2").Value = 1.95
$ws.Range("T2").Value = 1.95
$ws.Range("U2").Value = 2.38
$ws.Range("V2").Value = 1.57
$ws.Range("W2").Value = 1.29
$ws.Range("X2").Value = 3.5
$ws.Range("AF2").Value = 21
$ws.Range("AK2").Value = 126
$ws.Range("G3").Value = 2.55
$ws.Range("I3").Value = 2.63
$ws.Range("O3").Value = 1.29
$ws.Range("P3").Value = 3.75
$ws.Range("Q3").Value = 1.89
$ws.Range("R3").Value = 2.01
$ws.Range("W3").Value = 1.36
$ws.Range("X3").Value = 3
$ws.Range("Y3").Value = 1.67
$ws.Range("Z3").Value = 2.1
$ws.Range("AG3").Value = 12
$ws.Range("AL3").Value = 10
$ws.Range("AQ3").Value = 29
$ws.Range("G4").Value = 2.63
$ws.Range("I4").Value = 2.6
$ws.Range("AG4").Value = 10
$ws.Range("G5").Value = 1.17
$ws.Range("K5").Value = 3.25
$ws.Range("AC5").Value = 12
$ws.Range("AM5").Value = 81
$ws.Range("H6").Value = 3.5
$ws.Range("J6").Value = 2.55
$ws.Range("N6").Value = 7.9
$ws.Range("O6").Value = 1.25
$ws.Range("P6").Value = 3.6
$ws.Range("R6").Value = 1.98
$ws.Range("U6").Value = 2.77
$ws.Range("V6").Value = 1.39
$ws.Range("W6").Value = 1.37
$ws.Range("X6").Value = 2.85
$ws.Range("AA6").Value = 8.25
$ws.Range("AB6").Value = 10
$ws.Range("AF6").Value = 23
$ws.Range("AG6").Value = 7.9
$ws.Range("AI6").Value = 13
$ws.Range("AK6").Value = 350
$ws.Range("G7").Value = 1.7
$ws.Range("M7").Value = 1.1
$ws.Range("N7").Value = 7
$ws.Range("S7").Value = 4.1
$ws.Range("Y7").Value = 2.5
$ws.Range("Z7").Value = 1.5
$ws.Range("AE7").Value = 19
$ws.Range("AG7").Value = 6
$ws.Range("AH7").Value = 6.5
$ws.Range("AR7").Value = 1.95
$ws.Range("AS7").Value = 1.9
$ws.Range("G8").Value = 2.63
$ws.Range("I8").Value = 2.8
$ws.Range("J8").Value = 3.25
$ws.Range("L8").Value = 3.4
$ws.Range("AA8").Value = 8.5
$ws.Range("AB8").Value = 13
$ws.Range("AL8").Value = 8.5
$ws.Range("G9").Value = 3.1
$ws.Range("H9").Value = 3.5
$ws.Range("J9").Value = 3.75
$ws.Range("O10").Value = 1.36
$ws.Range("P10").Value = 3
$ws.Range("Q10").Value = 2.1
$ws.Range("R10").Value = 1.7
$ws.Range("G12").Value = 2.88
$ws.Range("M12").Value = 1.11
$ws.Range("N12").Value = 6.5
$ws.Range("AC12").Value = 11
$ws.Range("AD12").Value = 29
$ws.Range("AG12").Value = 6.5
$ws.Range("AP12").Value = 26
$ws.Range("AR12").Value = 1.95
$ws.Range("AS12").Value = 1.9
$ws.Range("H13").Value = 3.75
$ws.Range("I13").Value = 3.8
$ws.Range("J13").Value = 2.38
$ws.Range("L13").Value = 4.33
$ws.Range("AD13").Value = 15
$ws.Range("AK13").Value = 151
$ws.Range("AM13").Value = 23
$ws.Range("G15").Value = 4.33
$ws.Range("I15").Value = 1.67
$ws.Range("J15").Value = 4.75
$ws.Range("L15").Value = 2.2
$ws.Range("M15").Value = 1.03
$ws.Range("N15").Value = 10
$ws.Range("O15").Value = 1.2
$ws.Range("P15").Value = 4.33
$ws.Range("Q15").Value = 1.67
$ws.Range("R15").Value = 2.15
$ws.Range("U15").Value = 2.63
$ws.Range("V15").Value = 1.44
$ws.Range("AA15").Value = 15
$ws.Range("AD15").Value = 51
$ws.Range("AF15").Value = 41
$ws.Range("AH15").Value = 8
$ws.Range("AJ15").Value = 51
$ws.Range("H16").Value = 9.5
$ws.Range("I16").Value = 12
$ws.Range("J16").Value = 1.4
$ws.Range("K16").Value = 3.25
$ws.Range("N16").Value = 17
$ws.Range("Y16").Value = 2
$ws.Range("Z16").Value = 1.73
$ws.Range("AB16").Value = 8
$ws.Range("AD16").Value = 7.5
$ws.Range("AF16").Value = 29
$ws.Range("AL16").Value = 34
$ws.Range("G18").Value = 2.15
$ws.Range("AA18").Value = 8.5
$ws.Range("AH18").Value = 7
$ws.Range("AO18").Value = 34
$ws.Range("G19").Value = 1.14
$ws.Range("H19").Value = 9
$ws.Range("I19").Value = 13
$ws.Range("J19").Value = 1.44
$ws.Range("K19").Value = 3.25
$ws.Range("L19").Value = 11
$ws.Range("N19").Value = 29
$ws.Range("O19").Value = 1.08
$ws.Range("P19").Value = 8
$ws.Range("AG19").Value = 29
$ws.Range("AI19").Value = 29
$ws.Range("AO19").Value = 201
$ws.Range("G21").Value = 2.8
$ws.Range("I21").Value = 2.55
$ws.Range("M21").Value = 1.07
$ws.Range("N21").Value = 7.5
$ws.Range("O21").Value = 1.33
$ws.Range("P21").Value = 3.25
$ws.Range("AA21").Value = 8.5
$ws.Range("G22").Value = 1.44
$ws.Range("J22").Value = 1.95
$ws.Range("Q22").Value = 2
$ws.Range("R22").Value = 1.85
$ws.Range("U22").Value = 3.4
$ws.Range("V22").Value = 1.3
$ws.Range("AD22").Value = 9
$ws.Range("AE22").Value = 13
$ws.Range("AL22").Value = 15
$ws.Range("J23").Value = 4.15
$ws.Range("N23").Value = 8.5
$ws.Range("O23").Value = 1.25
$ws.Range("P23").Value = 3.7
$ws.Range("Q23").Value = 1.75
$ws.Range("R23").Value = 2
$ws.Range("U23").Value = 2.8
$ws.Range("V23").Value = 1.4
$ws.Range("W23").Value = 1.37
$ws.Range("X23").Value = 2.95
$ws.Range("Y23").Value = 1.7
$ws.Range("Z23").Value = 2.05
$ws.Range("AG23").Value = 8.5
$ws.Range("AJ23").Value = 65
$ws.Range("AK23").Value = 500
$ws.Range("AL23").Value = 7.8
$ws.Range("AO23").Value = 16.5
$ws.Range("AQ23").Value = 26
